$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Steps")
$ws.Range("A2").Value = 10000

$ws.Range("A3").Value = 12000
$ws.Range("B3").Value = 12000
$ws.Range("C3").Value = 2.2

$ws.Range("A4").Value = 15000
$ws.Range("B4").Value = 20000
$ws.Range("C4").Value = 3

$ws.Range("A5").Value = 20000
$ws.Range("B5").Value = 20000
$ws.Range("C5").Value = 4

$ws.Range("A6").Value = 13000
$ws.Range("B6").Value = 15000
$ws.Range("B6").HorizontalAlignment = -4108
$ws.Range("C6").Value = 2.5

$ws.Range("A7").Value = 14000
$ws.Range("B7").Value = 15000
$ws.Range("C7").Value = 3

$ws.Range("A8").Value = 15000
$ws.Range("B8").Value = 15000
$ws.Range("C8").Value = 3

$ws.Range("A9").Value = 21000
$ws.Range("B9").Value = 20000
$ws.Range("C9").Value = 3.5

$ws.Range("A10").Value = 18000
$ws.Range("B10").Value = 20000
$ws.Range("C10").Value = 3

$ws.Range("A11").Value = 15000
$ws.Range("B11").Value = 15000
$ws.Range("C11").Value = 2

$ws.Range("C11").Select()

$ws.Activate()

$homeWs = $wb.Worksheets.Item("Home")
$homeWs.Range("A2:A11").Select()
